$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C26").Value = 20
$ws.Range("E26").Value = 1505
